$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Update the CLIENTES list for the "RECURRENTE_<_100K" group (row 5) with
# the newly added client accounts.
$ws.Range("B5").Value = "12.20.265.938.947.970.999.1000.20094.20095.20108.20117.20140.20142.20156.20160.20169.20179.20242.20244.20249.20294.20344.20362.20364.20366.70106.60109.20189.90622.60286.60284.60281.60170.60241.60250.60282.60242.90631.60266.60283.20397.1035.1031.50611.10166.60290.60289.60291.60287.00989.1034.50832.60285.50831.951.1038.1037.50445.1023.1030.1033.1036.1029.50630.1024.1020.1027.1025.1028.50830.1021.1026.1032.20183.1036.1039.1040"

# Update the active selection to B6 as left by the author after the edit.
$ws.Range("B6").Select()
